$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-06 03:50:32"
$wsZhCn.Range("G4").Value = "2016-02-06 03:51:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-06 03:50:46"
$wsDeDe.Range("G4").Value = "2016-02-06 03:51:38"
